# 9th Stab - Cosmetic Changes
# Insert two new "week" columns (Jun_17, Jun_15) ahead of the existing
# Jun_13 / Jun_10 columns, shifting the older data to the right, and
# record the new 6/15/2018 upgrade for BidaskClub.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the highlight fill colour used for "new rating this week" cells
# (e.g. B11 currently holds the most-recent-week highlighted rating).
$highlightColor = $ws.Range("B11").Interior.Color()

# Insert two blank columns before column B. This shifts the existing
# "Jun_13" column (B) to D and the existing "Jun_10" column (C) to E,
# carrying their values and formatting (including the highlighted cells)
# along with them.
$ws.Columns("B:C").Insert()

# New header row for the two freshly inserted week columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the new columns with the default "UN" (unchanged) rating used
# throughout the sheet.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Record this week's (6/15/2018) new rating for BidaskClub (row 22) and
# highlight it the same way the other "new this week" cells are highlighted.
$ws.Range("C22").Value = "6/15/2018,Upgrades,Sell -> Hold,"
$ws.Range("C22").Interior.Color = $highlightColor

# Match the column widths used across the week columns.
$ws.Columns("C:E").ColumnWidth = 7.140625
